$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "RG-20002"
$ws.Range("B2").Value = "P-388412033222"

$ws.Range("B3").Select()
